# "finestra incidenza 7gg centrata su ultimo g"
#
# The rolling-7-day-sum columns C ("somma mobile 7gg.") and D ("somma
# mobile 7gg. per 100mila abitanti") were computed with the window
# trailing the day in column A/B. This change re-centers the 7-day
# window on the last day, which (for a fixed-width trailing window)
# is equivalent to shifting every C/D value down by three rows: the
# value that used to land on row r now lands on row r+3.
#
# We therefore rebuild C/D for rows 5..184 from the OLD values that
# used to sit three rows above (rows 2..181), walking bottom-to-top so
# we never read a cell we've already overwritten. Rows 2..4 (and the
# old rows 2..4 that feed rows 5..7) were blank, so the first three
# destination rows (5,6,7) become blank too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 5
$lastDataRow = 184
$shift = 3

for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $srcRow = $r - $shift

    $srcC = $ws.Cells.Item($srcRow, 3)
    $srcD = $ws.Cells.Item($srcRow, 4)
    $dstC = $ws.Cells.Item($r, 3)
    $dstD = $ws.Cells.Item($r, 4)

    $srcCText = $srcC.Text
    $srcDText = $srcD.Text

    if ($srcCText.Length -eq 0) {
        # Source was blank -> destination becomes a blank text cell too
        # (matches the existing blank cells such as C2:D4), rather than
        # a cleared/empty cell.
        $dstC.Value = "'"
        $dstC.Style = "Normal"
    } else {
        $dstC.Value2 = $srcC.Value2
    }

    if ($srcDText.Length -eq 0) {
        $dstD.Value = "'"
        $dstD.Style = "Normal"
    } else {
        $dstD.Value2 = $srcD.Value2
    }
}
